# Slide 8: content placeholder "Espace réservé du contenu 2" gets a new
# paragraph ("It will compress data as much as he can") inserted right
# after the "So do SAI need to" paragraph.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(8)

# Find the shape that contains the "So do SAI ... need to" text.
$targetShape = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $shp = $s.Shapes.Item($i)
    if ($shp.HasTextFrame -and $shp.TextFrame.HasText -and $shp.TextFrame.TextRange.Text -like "*So do SAI*need*to*") {
        $targetShape = $shp
    }
}

$tr = $targetShape.TextFrame.TextRange

# Find the paragraph that reads "So do SAI need to".
$targetPara = $null
for ($i = 1; $i -le $tr.Paragraphs().Count; $i++) {
    $para = $tr.Paragraphs($i)
    if ($para.Text -like "*So do SAI*need*to*") {
        $targetPara = $para
    }
}

# Typing a new line right after "So do SAI need to": start a new
# paragraph and add the sentence.
[void]$targetPara.InsertAfter("`rIt will compress data as much as he can")
